$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking values (e.g. "607.67")
# are stored as literal text, matching the source data (inlineStr cells), not
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '66.286.50'
$ws.Range('E2').Value = '  +1.28%  '
$ws.Range('D3').Value = '3.573.40'
$ws.Range('E3').Value = '  +5.38%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '607.67'
$ws.Range('E5').Value = '  +2.19%  '
$ws.Range('D6').Value = '145.47'
$ws.Range('E6').Value = '  +2.62%  '
$ws.Range('D7').Value = '3.568.35'
$ws.Range('E7').Value = '  +5.26%  '
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').Value = '0.487'
$ws.Range('E9').Value = '  +3.83%  '
$ws.Range('E10').Value = '  +2.32%  '
$ws.Range('D11').Value = '8.04'
$ws.Range('E11').Value = '  +1.70%  '
$ws.Range('D12').Value = '0.413'
$ws.Range('E12').Value = '  +1.62%  '
$ws.Range('D13').Value = '4.172.62'
$ws.Range('E13').Value = '  +5.15%  '
$ws.Range('D14').Value = '0.0000209'
$ws.Range('E14').Value = '  +4.77%  '
$ws.Range('D15').Value = '30.21'
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('D16').Value = '3.549.72'
$ws.Range('E16').Value = '  +4.57%  '
$ws.Range('D17').Value = '66.388.53'
$ws.Range('E17').Value = '  +1.37%  '
$ws.Range('E18').Value = '  -0.71%  '
$ws.Range('D19').Value = '11.53'
$ws.Range('E19').Value = '  +11.56%  '
$ws.Range('D20').Value = '6.23'
$ws.Range('E20').Value = '  +2.05%  '
$ws.Range('D21').Value = '14.98'
$ws.Range('E21').Value = '  +1.64%  '
$ws.Range('D22').Value = '431.83'
$ws.Range('E22').Value = '  +3.91%  '
$ws.Range('D23').Value = '0.611'
$ws.Range('E23').Value = '  +5.44%  '
$ws.Range('D24').Value = '78.73'
$ws.Range('E24').Value = '  +1.77%  '
$ws.Range('D25').Value = '3.709.38'
$ws.Range('E25').Value = '  +5.19%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  +9.00%  '
$ws.Range('E28').Value = '  +4.44%  '
$ws.Range('D29').Value = '8.05'
$ws.Range('E29').Value = '  +2.91%  '
$ws.Range('D30').Value = '9.18'
$ws.Range('E30').Value = '  -0.83%  '
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('E32').Value = '  +1.30%  '
$ws.Range('D33').Value = '0.159'
$ws.Range('E33').Value = '  -0.51%  '
$ws.Range('D34').Value = '3.561.40'
$ws.Range('E34').Value = '  +5.06%  '
$ws.Range('D35').Value = '25.48'
$ws.Range('E35').Value = '  +4.62%  '
$ws.Range('E36').Value = '  +4.60%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').Value = '7.92'
$ws.Range('E38').Value = '  +4.91%  '
$ws.Range('D39').Value = '5.66'
$ws.Range('E39').Value = '  +2.36%  '
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('D41').Value = '171.25'
$ws.Range('E41').Value = '  +0.94%  '
$ws.Range('D42').Value = '0.0858'
$ws.Range('E42').Value = '  +0.27%  '
$ws.Range('E43').Value = '  +3.62%  '
$ws.Range('D44').Value = '0.899'
$ws.Range('E44').Value = '  +3.49%  '
$ws.Range('D45').Value = '1.95'
$ws.Range('E45').Value = '  +1.85%  '
$ws.Range('E46').Value = '  +1.47%  '
$ws.Range('E47').Value = '  +4.52%  '
$ws.Range('D48').Value = '26.03'
$ws.Range('E48').Value = '  -2.20%  '
$ws.Range('E49').Value = '  +4.92%  '
$ws.Range('E50').Value = '  +1.27%  '
$ws.Range('D51').Value = '23.36'
$ws.Range('E51').Value = '  +15.79%  '

# Restore the original (default/"Normal") cell style now that the text values
# are committed, so no residual explicit style is left on the data cells.
$ws.Range("D2:D51").Style = "Normal"
